$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = 'Sending cluster'
$ws.Cells.Item(1,2).Value = 'Ligand symbol'
$ws.Cells.Item(1,3).Value = 'Receptor symbol'
$ws.Cells.Item(1,4).Value = 'Target cluster'
$ws.Cells.Item(1,5).Value = 'Ligand-expressing cells'
$ws.Cells.Item(1,6).Value = 'Ligand detection rate'
$ws.Cells.Item(1,7).Value = 'Ligand average expression value'
$ws.Cells.Item(1,8).Value = 'Ligand total expression value'
$ws.Cells.Item(1,9).Value = 'Ligand derived specificity of average expression value'
$ws.Cells.Item(1,10).Value = 'Ligand derived specificity of total expression value'
$ws.Cells.Item(1,11).Value = 'Receptor-expressing cells'
$ws.Cells.Item(1,12).Value = 'Receptor detection rate'
$ws.Cells.Item(1,13).Value = 'Receptor average expression value'
$ws.Cells.Item(1,14).Value = 'Receptor total expression value'
$ws.Cells.Item(1,15).Value = 'Receptor derived specificity of average expression value'
$ws.Cells.Item(1,16).Value = 'Receptor derived specificity of total expression value'
$ws.Cells.Item(1,17).Value = 'Edge average expression weight'
$ws.Cells.Item(1,18).Value = 'Edge total expression weight'
$ws.Cells.Item(1,19).Value = 'Edge average expression derived specificity'
$ws.Cells.Item(1,20).Value = 'Edge total expression derived specificity'

# Row 2
$ws.Cells.Item(2,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(2,2).Value = 'Fgf13'
$ws.Cells.Item(2,3).Value = 'Scn8a'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.2647293333333333
$ws.Cells.Item(2,8).Value = 0.794188
$ws.Cells.Item(2,9).Value = 0.1169441137252306
$ws.Cells.Item(2,10).Value = 0.1580772484350275
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.1169926666666667
$ws.Cells.Item(2,14).Value = 0.350978
$ws.Cells.Item(2,15).Value = 0.1389115695260365
$ws.Cells.Item(2,16).Value = 0.1396447152103714
$ws.Cells.Item(2,17).Value = 0.03097139065155556
$ws.Cells.Item(2,18).Value = 0.278742515864
$ws.Cells.Item(2,19).Value = 0.0162448903844031
$ws.Cells.Item(2,20).Value = 0.02207465233894856

# Row 3
$ws.Cells.Item(3,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(3,2).Value = 'Fgf13'
$ws.Cells.Item(3,3).Value = 'Scn8a'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.2647293333333333
$ws.Cells.Item(3,8).Value = 0.794188
$ws.Cells.Item(3,9).Value = 0.1169441137252306
$ws.Cells.Item(3,10).Value = 0.1580772484350275
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.679503
$ws.Cells.Item(3,14).Value = 2.038509
$ws.Cells.Item(3,15).Value = 0.8068097848952103
$ws.Cells.Item(3,16).Value = 0.8110679551390089
$ws.Cells.Item(3,17).Value = 0.179884376188
$ws.Cells.Item(3,18).Value = 1.618959385692
$ws.Cells.Item(3,19).Value = 0.09435165523941434
$ws.Cells.Item(3,20).Value = 0.1282113906421989

# Row 4
$ws.Cells.Item(4,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4,2).Value = 'Fgf13'
$ws.Cells.Item(4,3).Value = 'Scn8a'
$ws.Cells.Item(4,4).Value = 'MuSCs'
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.2647293333333333
$ws.Cells.Item(4,8).Value = 0.794188
$ws.Cells.Item(4,9).Value = 0.1169441137252306
$ws.Cells.Item(4,10).Value = 0.1580772484350275
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.5
$ws.Cells.Item(4,13).Value = 0.013265
$ws.Cells.Item(4,14).Value = 0.02653
$ws.Cells.Item(4,15).Value = 0.01575023479901482
$ws.Cells.Item(4,16).Value = 0.01055557412296826
$ws.Cells.Item(4,17).Value = 0.003511634606666667
$ws.Cells.Item(4,18).Value = 0.02106980764
$ws.Cells.Item(4,19).Value = 0.001841897249535074
$ws.Cells.Item(4,20).Value = 0.001668596113010802

# Row 5
$ws.Cells.Item(5,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(5,2).Value = 'Fgf13'
$ws.Cells.Item(5,3).Value = 'Scn8a'
$ws.Cells.Item(5,4).Value = 'Resolving-Mac'
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.2647293333333333
$ws.Cells.Item(5,8).Value = 0.794188
$ws.Cells.Item(5,9).Value = 0.1169441137252306
$ws.Cells.Item(5,10).Value = 0.1580772484350275
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.032449
$ws.Cells.Item(5,14).Value = 0.097347
$ws.Cells.Item(5,15).Value = 0.03852841077973854
$ws.Cells.Item(5,16).Value = 0.03873175552765139
$ws.Cells.Item(5,17).Value = 0.008590202137333332
$ws.Cells.Item(5,18).Value = 0.077311819236
$ws.Cells.Item(5,19).Value = 0.004505670851878147
$ws.Cells.Item(5,20).Value = 0.0061226093408693

# Row 6
$ws.Cells.Item(6,1).Value = 'MuSCs'
$ws.Cells.Item(6,2).Value = 'Fgf13'
$ws.Cells.Item(6,3).Value = 'Scn8a'
$ws.Cells.Item(6,4).Value = 'ECs'
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.7671255
$ws.Cells.Item(6,8).Value = 3.534251
$ws.Cells.Item(6,9).Value = 0.7806272272009463
$ws.Cells.Item(6,10).Value = 0.7034665260098926
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.1169926666666667
$ws.Cells.Item(6,14).Value = 0.350978
$ws.Cells.Item(6,15).Value = 0.1389115695260365
$ws.Cells.Item(6,16).Value = 0.1396447152103714
$ws.Cells.Item(6,17).Value = 0.2067407245796667
$ws.Cells.Item(6,18).Value = 1.240444347478
$ws.Cells.Item(6,19).Value = 0.1084381533452414
$ws.Cells.Item(6,20).Value = 0.09823538268468081

# Row 7
$ws.Cells.Item(7,1).Value = 'MuSCs'
$ws.Cells.Item(7,2).Value = 'Fgf13'
$ws.Cells.Item(7,3).Value = 'Scn8a'
$ws.Cells.Item(7,4).Value = 'FAPs'
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.7671255
$ws.Cells.Item(7,8).Value = 3.534251
$ws.Cells.Item(7,9).Value = 0.7806272272009463
$ws.Cells.Item(7,10).Value = 0.7034665260098926
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.679503
$ws.Cells.Item(7,14).Value = 2.038509
$ws.Cells.Item(7,15).Value = 0.8068097848952103
$ws.Cells.Item(7,16).Value = 0.8110679551390089
$ws.Cells.Item(7,17).Value = 1.2007670786265
$ws.Cells.Item(7,18).Value = 7.204602471758999
$ws.Cells.Item(7,19).Value = 0.62981768526134
$ws.Cells.Item(7,20).Value = 0.5705591567595859

# Row 8
$ws.Cells.Item(8,1).Value = 'MuSCs'
$ws.Cells.Item(8,2).Value = 'Fgf13'
$ws.Cells.Item(8,3).Value = 'Scn8a'
$ws.Cells.Item(8,4).Value = 'MuSCs'
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.7671255
$ws.Cells.Item(8,8).Value = 3.534251
$ws.Cells.Item(8,9).Value = 0.7806272272009463
$ws.Cells.Item(8,10).Value = 0.7034665260098926
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.5
$ws.Cells.Item(8,13).Value = 0.013265
$ws.Cells.Item(8,14).Value = 0.02653
$ws.Cells.Item(8,15).Value = 0.01575023479901482
$ws.Cells.Item(8,16).Value = 0.01055557412296826
$ws.Cells.Item(8,17).Value = 0.0234409197575
$ws.Cells.Item(8,18).Value = 0.09376367903000001
$ws.Cells.Item(8,19).Value = 0.01229506211891879
$ws.Cells.Item(8,20).Value = 0.007425493058324401

# Row 9
$ws.Cells.Item(9,1).Value = 'MuSCs'
$ws.Cells.Item(9,2).Value = 'Fgf13'
$ws.Cells.Item(9,3).Value = 'Scn8a'
$ws.Cells.Item(9,4).Value = 'Resolving-Mac'
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.7671255
$ws.Cells.Item(9,8).Value = 3.534251
$ws.Cells.Item(9,9).Value = 0.7806272272009463
$ws.Cells.Item(9,10).Value = 0.7034665260098926
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.032449
$ws.Cells.Item(9,14).Value = 0.097347
$ws.Cells.Item(9,15).Value = 0.03852841077973854
$ws.Cells.Item(9,16).Value = 0.03873175552765139
$ws.Cells.Item(9,17).Value = 0.0573414553495
$ws.Cells.Item(9,18).Value = 0.344048732097
$ws.Cells.Item(9,19).Value = 0.03007632647544635
$ws.Cells.Item(9,20).Value = 0.02724649350730137

# Row 10
$ws.Cells.Item(10,1).Value = 'Resolving-Mac'
$ws.Cells.Item(10,2).Value = 'Fgf13'
$ws.Cells.Item(10,3).Value = 'Scn8a'
$ws.Cells.Item(10,4).Value = 'ECs'
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.2318703333333333
$ws.Cells.Item(10,8).Value = 0.695611
$ws.Cells.Item(10,9).Value = 0.1024286590738231
$ws.Cells.Item(10,10).Value = 0.1384562255550801
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.1169926666666667
$ws.Cells.Item(10,14).Value = 0.350978
$ws.Cells.Item(10,15).Value = 0.1389115695260365
$ws.Cells.Item(10,16).Value = 0.1396447152103714
$ws.Cells.Item(10,17).Value = 0.02712712861755556
$ws.Cells.Item(10,18).Value = 0.244144157558
$ws.Cells.Item(10,19).Value = 0.01422852579639207
$ws.Cells.Item(10,20).Value = 0.01933468018674211

# Row 11
$ws.Cells.Item(11,1).Value = 'Resolving-Mac'
$ws.Cells.Item(11,2).Value = 'Fgf13'
$ws.Cells.Item(11,3).Value = 'Scn8a'
$ws.Cells.Item(11,4).Value = 'FAPs'
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.2318703333333333
$ws.Cells.Item(11,8).Value = 0.695611
$ws.Cells.Item(11,9).Value = 0.1024286590738231
$ws.Cells.Item(11,10).Value = 0.1384562255550801
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.679503
$ws.Cells.Item(11,14).Value = 2.038509
$ws.Cells.Item(11,15).Value = 0.8068097848952103
$ws.Cells.Item(11,16).Value = 0.8110679551390089
$ws.Cells.Item(11,17).Value = 0.157556587111
$ws.Cells.Item(11,18).Value = 1.418009283999
$ws.Cells.Item(11,19).Value = 0.08264044439445603
$ws.Cells.Item(11,20).Value = 0.1122974077372242

# Row 12
$ws.Cells.Item(12,1).Value = 'Resolving-Mac'
$ws.Cells.Item(12,2).Value = 'Fgf13'
$ws.Cells.Item(12,3).Value = 'Scn8a'
$ws.Cells.Item(12,4).Value = 'MuSCs'
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.2318703333333333
$ws.Cells.Item(12,8).Value = 0.695611
$ws.Cells.Item(12,9).Value = 0.1024286590738231
$ws.Cells.Item(12,10).Value = 0.1384562255550801
$ws.Cells.Item(12,11).Value = 1
$ws.Cells.Item(12,12).Value = 0.5
$ws.Cells.Item(12,13).Value = 0.013265
$ws.Cells.Item(12,14).Value = 0.02653
$ws.Cells.Item(12,15).Value = 0.01575023479901482
$ws.Cells.Item(12,16).Value = 0.01055557412296826
$ws.Cells.Item(12,17).Value = 0.003075759971666667
$ws.Cells.Item(12,18).Value = 0.01845455983
$ws.Cells.Item(12,19).Value = 0.001613275430560953
$ws.Cells.Item(12,20).Value = 0.00146148495163306

# Row 13
$ws.Cells.Item(13,1).Value = 'Resolving-Mac'
$ws.Cells.Item(13,2).Value = 'Fgf13'
$ws.Cells.Item(13,3).Value = 'Scn8a'
$ws.Cells.Item(13,4).Value = 'Resolving-Mac'
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.2318703333333333
$ws.Cells.Item(13,8).Value = 0.695611
$ws.Cells.Item(13,9).Value = 0.1024286590738231
$ws.Cells.Item(13,10).Value = 0.1384562255550801
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.032449
$ws.Cells.Item(13,14).Value = 0.097347
$ws.Cells.Item(13,15).Value = 0.03852841077973854
$ws.Cells.Item(13,16).Value = 0.03873175552765139
$ws.Cells.Item(13,17).Value = 0.007523960446333333
$ws.Cells.Item(13,18).Value = 0.067715644017
$ws.Cells.Item(13,19).Value = 0.00394641345241405
$ws.Cells.Item(13,20).Value = 0.00536265267948072

